# Replace computed ratio formulas on the PoFDCtAE sheet with a constant
# value of 1 for the diagonal (From type == To type) cells that previously
# pulled their ratio from the "Data from BFPIaE" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PoFDCtAE")

# Cell -> replace formula with the literal value 1
$cells = @("C3", "D4", "I9", "J10", "K11", "L12", "M13", "N14", "S19", "T20")

foreach ($addr in $cells) {
    $ws.Range($addr).Value = 1
}

# Keep the saved selection state consistent with the source workbook
# (active cell on the "PoFDCtAE" sheet's bottom-right pane ends up on U20),
# then restore the originally active sheet ("About") so the tab selection
# in the saved file is unchanged.
$originalActive = $wb.ActiveSheet
$ws.Activate()
$ws.Range("U20").Select()
$originalActive.Activate()
